$d = $word.ActiveDocument

# Locate the paragraph that holds "This is a Microsoft word document." by
# searching the document content; this keeps the script robust to the
# target text's exact position. Find.Execute collapses $searchRange to the
# matched text when it succeeds.
$needle = "This is a Microsoft word document."
$searchRange = $d.Content
$found = $searchRange.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target paragraph text: $needle"
}

# Grab the paragraph that contains the matched text and compute its range
# without the trailing paragraph mark, so we can rebuild its run content.
$para = $searchRange.Paragraphs(1)
$paraStart = $para.Range.Start
$paraEnd = $para.Range.End - 1
$paraRange = $d.Range($paraStart, $paraEnd)

# New runs to append after the existing sentence: " (" / "Changed main" / ")"
# expressed as separate <w:r> elements, matching the target OOXML diff.
$openingXml = '<w:r><w:t xml:space="preserve"> (</w:t></w:r>'
$middleXml  = '<w:r><w:t>Changed main</w:t></w:r>'
$closingXml = '<w:r><w:t>)</w:t></w:r>'

$existingRunXml = '<w:r><w:t>' + $paraRange.Text + '</w:t></w:r>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' + $existingRunXml + $openingXml + $middleXml + $closingXml + '</w:p></w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$paraRange.InsertXML($packageXml)
